$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the existing contents; we'll rewrite the whole table with the new
# layout (two new columns "ownTeam"/"oppTeam" inserted after "result", all
# rows reordered, and five additional match rows appended).
$ws.Cells.Clear()

$rows = @(
    ,@("venue","date","result","ownTeam","oppTeam","batsman","totalRuns","totalBalls","total4s","total6s","sr")
    ,@(" Dubai (DSC)"," October 27 2020","Sunrisers won by 88 runs","Delhi Capitals","Sunrisers Hyderabad","Ajinkya Rahane ","26","19","3","1","136.84")
    ,@(" Dubai (DSC)"," November 05 2020","Mumbai won by 57 runs","Delhi Capitals","Mumbai Indians","Ajinkya Rahane ","0","3","0","0","0.00")
    ,@(" Dubai (DSC)"," November 10 2020","Mumbai won by 5 wickets (with 8 balls remaining)","Delhi Capitals","Mumbai Indians","Ajinkya Rahane ","2","4","0","0","50.00")
    ,@(" Dubai (DSC)"," October 14 2020","Capitals won by 13 runs","Delhi Capitals","Rajasthan Royals","Ajinkya Rahane ","2","9","0","0","22.22")
    ,@(" Sharjah"," October 17 2020","Capitals won by 5 wickets (with 1 ball remaining)","Delhi Capitals","Chennai Super Kings","Ajinkya Rahane ","8","10","1","0","80.00")
    ,@(" Abu Dhabi"," November 02 2020","Capitals won by 6 wickets (with 6 balls remaining)","Delhi Capitals","Royal Challengers Bangalore","Ajinkya Rahane ","60","46","5","1","130.43")
    ,@(" Abu Dhabi"," October 24 2020","KKR won by 59 runs","Delhi Capitals","Kolkata Knight Riders","Ajinkya Rahane ","0","1","0","0","0.00")
    ,@(" Abu Dhabi"," October 11 2020","Mumbai won by 5 wickets (with 2 balls remaining)","Delhi Capitals","Mumbai Indians","Ajinkya Rahane ","15","15","3","0","100.00")
)

for ($r = 0; $r -lt $rows.Count; $r++) {
    $rowValues = $rows[$r]
    for ($c = 0; $c -lt $rowValues.Count; $c++) {
        $cell = $ws.Cells.Item($r + 1, $c + 1)
        # Columns G..K (totalRuns, totalBalls, total4s, total6s, sr) hold
        # numeric-looking text in the source data; force text formatting so
        # they are written out as strings rather than numbers.
        if ($r -ge 1 -and $c -ge 6) {
            $cell.NumberFormat = "@"
        }
        $cell.Value = $rowValues[$c]
    }
}
